$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 8).Value = 1204.2222  # H2 (SET)
$ws.Cells.Item(2, 9).Value = 266.6  # I2 (SET)
$ws.Cells.Item(2, 10).Value = 2376.25  # J2 (SET)
$ws.Cells.Item(2, 11).Value = 266.6  # K2 (SET)
$ws.Cells.Item(2, 12).Value = 2376.25  # L2 (SET)
$ws.Cells.Item(2, 13).Value = -153.6  # M2 (SET)
$ws.Cells.Item(2, 14).Value = -2602.25  # N2 (SET)
$ws.Cells.Item(9, 8).Value = 1877.375  # H9 (SET)
$ws.Cells.Item(9, 10).Value = 849.6667  # J9 (SET)
$ws.Cells.Item(9, 12).Value = 849.6667  # L9 (SET)
$ws.Cells.Item(9, 14).Value = -1187.6667  # N9 (SET)
$ws.Cells.Item(32, 8).Value = 3572  # H32 (SET)
$ws.Cells.Item(32, 9).Value = 1187.25  # I32 (SET)
$ws.Cells.Item(32, 11).Value = 1187.25  # K32 (SET)
$ws.Cells.Item(32, 13).Value = -861.25  # M32 (SET)
$ws.Cells.Item(86, 8).Value = 177656.67  # H86 (SET)
$ws.Cells.Item(86, 9).Value = 13187.2  # I86 (SET)
$ws.Cells.Item(86, 11).Value = 13187.2  # K86 (SET)
$ws.Cells.Item(86, 13).Value = -12064.2  # M86 (SET)
$ws.Cells.Item(89, 8).Value = 177656.67  # H89 (SET)
$ws.Cells.Item(89, 9).Value = 13187.2  # I89 (SET)
$ws.Cells.Item(89, 11).Value = 65936  # K89 (SET)
$ws.Cells.Item(89, 13).Value = -60320  # M89 (SET)
$ws.Cells.Item(100, 8).Value = 3499.6667  # H100 (SET)
$ws.Cells.Item(100, 9).Value = 2500  # I100 (SET)
$ws.Cells.Item(100, 10).Value = 3999.5  # J100 (SET)
$ws.Cells.Item(100, 11).Value = 2500  # K100 (SET)
$ws.Cells.Item(100, 12).Value = 3999.5  # L100 (SET)
$ws.Cells.Item(100, 13).Value = -1959  # M100 (ADD)
$ws.Cells.Item(100, 14).Value = -5081.5  # N100 (SET)
$ws.Cells.Item(106, 8).Value = 1724.75  # H106 (SET)
$ws.Cells.Item(106, 9).Value = 1724.75  # I106 (SET)
$ws.Cells.Item(106, 11).Value = 1724.75  # K106 (SET)
$ws.Cells.Item(106, 13).Value = -1093.75  # M106 (SET)
$ws.Cells.Item(113, 8).Value = 3202.8  # H113 (SET)
$ws.Cells.Item(113, 10).Value = 3203.5  # J113 (SET)
$ws.Cells.Item(113, 12).Value = 3203.5  # L113 (SET)
$ws.Cells.Item(113, 14).Value = -9711.5  # N113 (SET)
$ws.Cells.Item(135, 8).Value = 325.46667  # H135 (SET)
$ws.Cells.Item(135, 9).Value = 325.46667  # I135 (SET)
$ws.Cells.Item(135, 11).Value = 2929.20003  # K135 (SET)
$ws.Cells.Item(135, 13).Value = -394.20003  # M135 (SET)
$ws.Cells.Item(137, 8).Value = 1658.3125  # H137 (SET)
$ws.Cells.Item(137, 9).Value = 1156.4286  # I137 (SET)
$ws.Cells.Item(137, 10).Value = 2048.6667  # J137 (SET)
$ws.Cells.Item(137, 11).Value = 3469.2858  # K137 (SET)
$ws.Cells.Item(137, 12).Value = 6146.000100000001  # L137 (SET)
$ws.Cells.Item(137, 13).Value = -919.2857999999997  # M137 (SET)
$ws.Cells.Item(137, 14).Value = -11246.0001  # N137 (SET)
$ws.Cells.Item(138, 8).Value = 5032.926  # H138 (SET)
$ws.Cells.Item(138, 9).Value = 2948.625  # I138 (SET)
$ws.Cells.Item(138, 10).Value = 5910.5264  # J138 (SET)
$ws.Cells.Item(138, 11).Value = 8845.875  # K138 (SET)
$ws.Cells.Item(138, 12).Value = 17731.5792  # L138 (SET)
$ws.Cells.Item(138, 13).Value = -3705.875  # M138 (SET)
$ws.Cells.Item(138, 14).Value = -28011.5792  # N138 (SET)
$ws.Cells.Item(141, 8).Value = 2089.682  # H141 (SET)
$ws.Cells.Item(141, 9).Value = 1862.2106  # I141 (SET)
$ws.Cells.Item(141, 11).Value = 5586.6318  # K141 (SET)
$ws.Cells.Item(141, 13).Value = -406.6318000000001  # M141 (SET)

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 2335782.8  # H32 (SET)
$ws.Cells.Item(32, 9).Value = 3336704.8  # I32 (SET)
$ws.Cells.Item(32, 10).Value = 584169.3  # J32 (SET)
$ws.Cells.Item(32, 11).Value = 3336704.8  # K32 (SET)
$ws.Cells.Item(32, 12).Value = 584169.3  # L32 (SET)
$ws.Cells.Item(32, 13).Value = -3336417.8  # M32 (SET)
$ws.Cells.Item(32, 14).Value = -584743.3  # N32 (SET)
$ws.Cells.Item(110, 8).Value = 6168012  # H110 (SET)
$ws.Cells.Item(110, 9).Value = 7401314.5  # I110 (SET)
$ws.Cells.Item(110, 10).Value = 1499  # J110 (SET)
$ws.Cells.Item(110, 11).Value = 7401314.5  # K110 (SET)
$ws.Cells.Item(110, 12).Value = 1499  # L110 (SET)
$ws.Cells.Item(110, 13).Value = -7399269.5  # M110 (SET)
$ws.Cells.Item(110, 14).Value = -5589  # N110 (SET)
$ws.Cells.Item(132, 8).Value = 1999.5  # H132 (SET)
$ws.Cells.Item(132, 9).Value = 1999.5  # I132 (SET)
$ws.Cells.Item(132, 11).Value = 5998.5  # K132 (SET)
$ws.Cells.Item(132, 13).Value = -3468.5  # M132 (SET)

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(20, 8).Value = 2914.8572  # H20 (SET)
$ws.Cells.Item(20, 9).Value = 2551.5  # I20 (SET)
$ws.Cells.Item(20, 10).Value = 3399.3333  # J20 (SET)
$ws.Cells.Item(20, 11).Value = 2551.5  # K20 (SET)
$ws.Cells.Item(20, 12).Value = 3399.3333  # L20 (SET)
$ws.Cells.Item(20, 13).Value = -2304.5  # M20 (SET)
$ws.Cells.Item(20, 14).Value = -3893.3333  # N20 (SET)
$ws.Cells.Item(22, 8).Value = 634  # H22 (SET)
$ws.Cells.Item(22, 9).Value = 200  # I22 (SET)
$ws.Cells.Item(22, 10).Value = 1502  # J22 (SET)
$ws.Cells.Item(22, 11).Value = 200  # K22 (SET)
$ws.Cells.Item(22, 12).Value = 1502  # L22 (SET)
$ws.Cells.Item(22, 13).Value = -27  # M22 (SET)
$ws.Cells.Item(22, 14).Value = -1848  # N22 (ADD)
$ws.Cells.Item(86, 8).Value = 1700.3334  # H86 (SET)
$ws.Cells.Item(86, 10).Value = 1650.8334  # J86 (SET)
$ws.Cells.Item(86, 12).Value = 1650.8334  # L86 (SET)
$ws.Cells.Item(86, 14).Value = -3896.8334  # N86 (SET)
$ws.Cells.Item(89, 8).Value = 1700.3334  # H89 (SET)
$ws.Cells.Item(89, 10).Value = 1650.8334  # J89 (SET)
$ws.Cells.Item(89, 12).Value = 8254.166999999999  # L89 (SET)
$ws.Cells.Item(89, 14).Value = -19486.167  # N89 (SET)
$ws.Cells.Item(105, 8).Value = 2288.7778  # H105 (SET)
$ws.Cells.Item(105, 9).Value = 2299.8572  # I105 (SET)
$ws.Cells.Item(105, 11).Value = 2299.8572  # K105 (SET)
$ws.Cells.Item(105, 13).Value = -552.8571999999999  # M105 (SET)
$ws.Cells.Item(107, 8).Value = 0  # H107 (SET)
$ws.Cells.Item(107, 9).Value = 0  # I107 (SET)
$ws.Cells.Item(107, 10).Value = 0  # J107 (SET)
$ws.Cells.Item(107, 11).Value = 0  # K107 (SET)
$ws.Cells.Item(107, 12).Value = 0  # L107 (SET)
$ws.Cells.Item(107, 13).ClearContents()  # M107
$ws.Cells.Item(107, 14).ClearContents()  # N107
$ws.Cells.Item(134, 8).Value = 2225  # H134 (SET)
$ws.Cells.Item(134, 9).Value = 2158.0715  # I134 (SET)
$ws.Cells.Item(134, 11).Value = 6474.2145  # K134 (SET)
$ws.Cells.Item(134, 13).Value = -3939.2145  # M134 (SET)

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(5, 8).Value = 889.2  # H5 (SET)
$ws.Cells.Item(5, 10).Value = 374  # J5 (SET)
$ws.Cells.Item(5, 12).Value = 374  # L5 (SET)
$ws.Cells.Item(5, 14).Value = -598  # N5 (SET)
$ws.Cells.Item(31, 8).Value = 2950  # H31 (SET)
$ws.Cells.Item(31, 9).Value = 0  # I31 (SET)
$ws.Cells.Item(31, 11).Value = 0  # K31 (SET)
$ws.Cells.Item(31, 13).ClearContents()  # M31
$ws.Cells.Item(34, 8).Value = 2950  # H34 (SET)
$ws.Cells.Item(34, 9).Value = 0  # I34 (SET)
$ws.Cells.Item(34, 11).Value = 0  # K34 (SET)
$ws.Cells.Item(34, 13).ClearContents()  # M34
$ws.Cells.Item(86, 8).Value = 12977.923  # H86 (SET)
$ws.Cells.Item(86, 10).Value = 11714.6  # J86 (SET)
$ws.Cells.Item(86, 12).Value = 11714.6  # L86 (SET)
$ws.Cells.Item(86, 14).Value = -13960.6  # N86 (SET)
$ws.Cells.Item(89, 8).Value = 12977.923  # H89 (SET)
$ws.Cells.Item(89, 10).Value = 11714.6  # J89 (SET)
$ws.Cells.Item(89, 12).Value = 58573  # L89 (SET)
$ws.Cells.Item(89, 14).Value = -69805  # N89 (SET)
$ws.Cells.Item(134, 8).Value = 3059.6  # H134 (SET)
$ws.Cells.Item(134, 9).Value = 2960.7144  # I134 (SET)
$ws.Cells.Item(134, 11).Value = 8882.143199999999  # K134 (SET)
$ws.Cells.Item(134, 13).Value = -6347.143199999999  # M134 (SET)

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(23, 8).Value = 443.85715  # H23 (SET)
$ws.Cells.Item(23, 10).Value = 591.1111  # J23 (SET)
$ws.Cells.Item(23, 12).Value = 1773.3333  # L23 (SET)
$ws.Cells.Item(23, 14).Value = -2243.3333  # N23 (SET)
$ws.Cells.Item(33, 8).Value = 546.6667  # H33 (SET)
$ws.Cells.Item(33, 9).Value = 270.5  # I33 (SET)
$ws.Cells.Item(33, 11).Value = 1623  # K33 (SET)
$ws.Cells.Item(33, 13).Value = -1340  # M33 (SET)
$ws.Cells.Item(81, 8).Value = 2200  # H81 (SET)
$ws.Cells.Item(81, 10).Value = 2400  # J81 (SET)
$ws.Cells.Item(81, 12).Value = 7200  # L81 (SET)
$ws.Cells.Item(81, 14).Value = -9446  # N81 (SET)
$ws.Cells.Item(84, 8).Value = 2200  # H84 (SET)
$ws.Cells.Item(84, 10).Value = 2400  # J84 (SET)
$ws.Cells.Item(84, 12).Value = 21600  # L84 (SET)
$ws.Cells.Item(84, 14).Value = -32832  # N84 (SET)
$ws.Cells.Item(132, 8).Value = 2330.8  # H132 (SET)
$ws.Cells.Item(132, 10).Value = 1832.7142  # J132 (SET)
$ws.Cells.Item(132, 12).Value = 16494.4278  # L132 (SET)
$ws.Cells.Item(132, 14).Value = -21554.4278  # N132 (SET)

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(70, 8).Value = 5302  # H70 (SET)
$ws.Cells.Item(70, 10).Value = 6000  # J70 (SET)
$ws.Cells.Item(70, 12).Value = 6000  # L70 (SET)
$ws.Cells.Item(70, 14).Value = -6540  # N70 (ADD)
$ws.Cells.Item(73, 8).Value = 5302  # H73 (SET)
$ws.Cells.Item(73, 10).Value = 6000  # J73 (SET)
$ws.Cells.Item(73, 12).Value = 6000  # L73 (SET)
$ws.Cells.Item(73, 14).Value = -7872  # N73 (ADD)
$ws.Cells.Item(80, 8).Value = 3210  # H80 (SET)
$ws.Cells.Item(80, 9).Value = 2760.6667  # I80 (SET)
$ws.Cells.Item(80, 10).Value = 3547  # J80 (SET)
$ws.Cells.Item(80, 11).Value = 2760.6667  # K80 (SET)
$ws.Cells.Item(80, 12).Value = 3547  # L80 (SET)
$ws.Cells.Item(80, 13).Value = -1762.6667  # M80 (SET)
$ws.Cells.Item(80, 14).Value = -5543  # N80 (SET)
$ws.Cells.Item(83, 8).Value = 3210  # H83 (SET)
$ws.Cells.Item(83, 9).Value = 2760.6667  # I83 (SET)
$ws.Cells.Item(83, 10).Value = 3547  # J83 (SET)
$ws.Cells.Item(83, 11).Value = 13803.3335  # K83 (SET)
$ws.Cells.Item(83, 12).Value = 17735  # L83 (SET)
$ws.Cells.Item(83, 13).Value = -8811.333500000001  # M83 (SET)
$ws.Cells.Item(83, 14).Value = -27719  # N83 (SET)
$ws.Cells.Item(132, 8).Value = 2623.3333  # H132 (SET)
$ws.Cells.Item(132, 9).Value = 2748.1  # I132 (SET)
$ws.Cells.Item(132, 10).Value = 1999.5  # J132 (SET)
$ws.Cells.Item(132, 11).Value = 8244.299999999999  # K132 (SET)
$ws.Cells.Item(132, 12).Value = 5998.5  # L132 (SET)
$ws.Cells.Item(132, 13).Value = -5714.299999999999  # M132 (SET)
$ws.Cells.Item(132, 14).Value = -11058.5  # N132 (ADD)

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(16, 8).Value = 5365  # H16 (SET)
$ws.Cells.Item(16, 9).Value = 456.25  # I16 (SET)
$ws.Cells.Item(16, 11).Value = 456.25  # K16 (SET)
$ws.Cells.Item(16, 13).Value = -286.25  # M16 (SET)
$ws.Cells.Item(22, 8).Value = 1868.4166  # H22 (SET)
$ws.Cells.Item(22, 9).Value = 1403.6666  # I22 (SET)
$ws.Cells.Item(22, 10).Value = 3262.6667  # J22 (SET)
$ws.Cells.Item(22, 11).Value = 1403.6666  # K22 (SET)
$ws.Cells.Item(22, 12).Value = 3262.6667  # L22 (SET)
$ws.Cells.Item(22, 13).Value = -1108.6666  # M22 (SET)
$ws.Cells.Item(22, 14).Value = -3852.6667  # N22 (SET)
$ws.Cells.Item(27, 8).Value = 1868.4166  # H27 (SET)
$ws.Cells.Item(27, 9).Value = 1403.6666  # I27 (SET)
$ws.Cells.Item(27, 10).Value = 3262.6667  # J27 (SET)
$ws.Cells.Item(27, 11).Value = 1403.6666  # K27 (SET)
$ws.Cells.Item(27, 12).Value = 3262.6667  # L27 (SET)
$ws.Cells.Item(27, 13).Value = -1296.6666  # M27 (SET)
$ws.Cells.Item(27, 14).Value = -3476.6667  # N27 (SET)
$ws.Cells.Item(46, 8).Value = 2955.6924  # H46 (SET)
$ws.Cells.Item(46, 9).Value = 2325.1667  # I46 (SET)
$ws.Cells.Item(46, 10).Value = 3496.1428  # J46 (SET)
$ws.Cells.Item(46, 11).Value = 2325.1667  # K46 (SET)
$ws.Cells.Item(46, 12).Value = 3496.1428  # L46 (SET)
$ws.Cells.Item(46, 13).Value = -2137.1667  # M46 (SET)
$ws.Cells.Item(46, 14).Value = -3872.1428  # N46 (SET)
$ws.Cells.Item(68, 8).Value = 3667  # H68 (SET)
$ws.Cells.Item(68, 9).Value = 3667  # I68 (SET)
$ws.Cells.Item(68, 11).Value = 3667  # K68 (SET)
$ws.Cells.Item(68, 13).Value = -2918  # M68 (SET)
$ws.Cells.Item(71, 8).Value = 3667  # H71 (SET)
$ws.Cells.Item(71, 9).Value = 3667  # I71 (SET)
$ws.Cells.Item(71, 11).Value = 18335  # K71 (SET)
$ws.Cells.Item(71, 13).Value = -14591  # M71 (SET)
$ws.Cells.Item(136, 8).Value = 4247  # H136 (SET)
$ws.Cells.Item(136, 9).Value = 3125  # I136 (SET)
$ws.Cells.Item(136, 11).Value = 9375  # K136 (SET)
$ws.Cells.Item(136, 13).Value = -6825  # M136 (SET)

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(55, 8).Value = 5101  # H55 (SET)
$ws.Cells.Item(55, 9).Value = 4682.6665  # I55 (SET)
$ws.Cells.Item(55, 11).Value = 4682.6665  # K55 (SET)
$ws.Cells.Item(55, 13).Value = -4405.6665  # M55 (SET)
$ws.Cells.Item(107, 8).Value = 517.5  # H107 (SET)
$ws.Cells.Item(107, 9).Value = 375.5  # I107 (SET)
$ws.Cells.Item(107, 10).Value = 801.5  # J107 (SET)
$ws.Cells.Item(107, 11).Value = 1126.5  # K107 (SET)
$ws.Cells.Item(107, 12).Value = 2404.5  # L107 (SET)
$ws.Cells.Item(107, 13).Value = 793.5  # M107 (ADD)
$ws.Cells.Item(107, 14).Value = -6244.5  # N107 (SET)
$ws.Cells.Item(132, 8).Value = 2348.5715  # H132 (SET)
$ws.Cells.Item(132, 9).Value = 2549.6  # I132 (SET)
$ws.Cells.Item(132, 10).Value = 1846  # J132 (SET)
$ws.Cells.Item(132, 11).Value = 7648.799999999999  # K132 (SET)
$ws.Cells.Item(132, 12).Value = 5538  # L132 (SET)
$ws.Cells.Item(132, 13).Value = -5118.799999999999  # M132 (SET)
$ws.Cells.Item(132, 14).Value = -10598  # N132 (SET)
